# Auto-generated Excel COM-interop script applying the Shinryu_Profits workbook update.
# For each affected leve-profit row (columns H-N), update the recalculated values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 200200240
$ws.Range("I18").Value = 250250140
$ws.Range("J18").Value = 600
$ws.Range("K18").Value = 250250140
$ws.Range("L18").Value = 600
$ws.Range("M18").Value = -250249856
$ws.Range("N18").Value = -1168

$ws.Range("H40").Value = 1669.6471
$ws.Range("I40").Value = 1783.5
$ws.Range("J40").Value = 1396.4
$ws.Range("K40").Value = 1783.5
$ws.Range("L40").Value = 1396.4
$ws.Range("M40").Value = -1608.5
$ws.Range("N40").Value = -1746.4

$ws.Range("H87").Value = 38772
$ws.Range("J87").Value = 47215
$ws.Range("L87").Value = 47215
$ws.Range("N87").Value = -49711

$ws.Range("H90").Value = 38772
$ws.Range("J90").Value = 47215
$ws.Range("L90").Value = 141645
$ws.Range("N90").Value = -154125

$ws.Range("H93").Value = 81998.75
$ws.Range("J93").Value = 81998.75
$ws.Range("L93").Value = 81998.75
$ws.Range("N93").Value = -86990.75

$ws.Range("H98").Value = 737.86664
$ws.Range("I98").Value = 465.66666
$ws.Range("J98").Value = 1826.6666
$ws.Range("K98").Value = 465.66666
$ws.Range("L98").Value = 1826.6666
$ws.Range("M98").Value = 1032.33334
$ws.Range("N98").Value = -4822.6666

$ws.Range("H99").Value = 284.83334
$ws.Range("I99").Value = 241.8
$ws.Range("J99").Value = 500
$ws.Range("K99").Value = 725.4000000000001
$ws.Range("L99").Value = 1500
$ws.Range("M99").Value = 772.5999999999999
$ws.Range("N99").Value = -4496

$ws.Range("H101").Value = 5490.9
$ws.Range("I101").Value = 392.83334
$ws.Range("J101").Value = 13138
$ws.Range("K101").Value = 1178.50002
$ws.Range("L101").Value = 39414
$ws.Range("M101").Value = 443.4999800000001
$ws.Range("N101").Value = -42658

$ws.Range("H122").Value = 737.86664
$ws.Range("I122").Value = 465.66666
$ws.Range("J122").Value = 1826.6666
$ws.Range("K122").Value = 1396.99998
$ws.Range("L122").Value = 5479.9998
$ws.Range("M122").Value = 1053.00002
$ws.Range("N122").Value = -10379.9998

$ws.Range("H137").Value = 1462.7028
$ws.Range("I137").Value = 1107.5358
$ws.Range("J137").Value = 2567.6667
$ws.Range("K137").Value = 3322.6074
$ws.Range("L137").Value = 7703.000100000001
$ws.Range("M137").Value = -772.6074000000003
$ws.Range("N137").Value = -12803.0001


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25645938
$ws.Range("I32").Value = 27030312
$ws.Range("K32").Value = 27030312
$ws.Range("M32").Value = -27030025

$ws.Range("H45").Value = 1740.5
$ws.Range("I45").Value = 1300
$ws.Range("J45").Value = 1810.0526
$ws.Range("K45").Value = 1300
$ws.Range("L45").Value = 1810.0526
$ws.Range("M45").Value = -923
$ws.Range("N45").Value = -2564.0526

$ws.Range("H61").Value = 1774.909
$ws.Range("I61").Value = 1452.4
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 1452.4
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -1240.4
$ws.Range("N61").Value = -5424

$ws.Range("H74").Value = 1518.9667
$ws.Range("I74").Value = 1473.125
$ws.Range("J74").Value = 1702.3334
$ws.Range("K74").Value = 1473.125
$ws.Range("L74").Value = 1702.3334
$ws.Range("M74").Value = -599.125
$ws.Range("N74").Value = -3450.3334

$ws.Range("H77").Value = 1518.9667
$ws.Range("I77").Value = 1473.125
$ws.Range("J77").Value = 1702.3334
$ws.Range("K77").Value = 7365.625
$ws.Range("L77").Value = 8511.666999999999
$ws.Range("M77").Value = -2997.625
$ws.Range("N77").Value = -17247.667

$ws.Range("H80").Value = 28850
$ws.Range("J80").Value = 28850
$ws.Range("L80").Value = 28850
$ws.Range("N80").Value = -30846

$ws.Range("H83").Value = 28850
$ws.Range("J83").Value = 28850
$ws.Range("L83").Value = 86550
$ws.Range("N83").Value = -96534

$ws.Range("H132").Value = 1717.3077
$ws.Range("I132").Value = 1384.1765
$ws.Range("J132").Value = 2346.5557
$ws.Range("K132").Value = 4152.529500000001
$ws.Range("L132").Value = 7039.6671
$ws.Range("M132").Value = -1622.529500000001
$ws.Range("N132").Value = -12099.6671

$ws.Range("H136").Value = 1774.909
$ws.Range("I136").Value = 1452.4
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 4357.200000000001
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -1807.200000000001
$ws.Range("N136").Value = -20100

$ws.Range("H139").Value = 50000
$ws.Range("J139").Value = 50000
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 10580
$ws.Range("J82").Value = 26750
$ws.Range("L82").Value = 26750
$ws.Range("N82").Value = -27516

$ws.Range("H85").Value = 10580
$ws.Range("J85").Value = 26750
$ws.Range("L85").Value = 26750
$ws.Range("N85").Value = -29402

$ws.Range("H99").Value = 1046.7878
$ws.Range("I99").Value = 651.9231
$ws.Range("J99").Value = 2513.4285
$ws.Range("K99").Value = 651.9231
$ws.Range("L99").Value = 2513.4285
$ws.Range("M99").Value = 846.0769
$ws.Range("N99").Value = -5509.4285

$ws.Range("H134").Value = 2000.2142
$ws.Range("I134").Value = 1916.4348
$ws.Range("K134").Value = 5749.3044
$ws.Range("M134").Value = -3214.3044


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2913.4092
$ws.Range("I16").Value = 2700.1
$ws.Range("J16").Value = 3091.1667
$ws.Range("K16").Value = 2700.1
$ws.Range("L16").Value = 3091.1667
$ws.Range("M16").Value = -2413.1
$ws.Range("N16").Value = -3665.1667

$ws.Range("H31").Value = 2338.3428
$ws.Range("I31").Value = 1846.6086
$ws.Range("J31").Value = 3280.8333
$ws.Range("K31").Value = 1846.6086
$ws.Range("L31").Value = 3280.8333
$ws.Range("M31").Value = -1551.6086
$ws.Range("N31").Value = -3870.8333

$ws.Range("H34").Value = 2338.3428
$ws.Range("I34").Value = 1846.6086
$ws.Range("J34").Value = 3280.8333
$ws.Range("K34").Value = 1846.6086
$ws.Range("L34").Value = 3280.8333
$ws.Range("M34").Value = -1644.6086
$ws.Range("N34").Value = -3684.8333

$ws.Range("H41").Value = 9032.091
$ws.Range("J41").Value = 9366.951999999999
$ws.Range("L41").Value = 9366.951999999999
$ws.Range("N41").Value = -10222.952

$ws.Range("H50").Value = 16200
$ws.Range("J50").Value = 17840
$ws.Range("L50").Value = 17840
$ws.Range("N50").Value = -19090

$ws.Range("H51").Value = 18500

$ws.Range("H60").Value = 9690.556
$ws.Range("I60").Value = 6350
$ws.Range("J60").Value = 10645
$ws.Range("K60").Value = 6350
$ws.Range("L60").Value = 10645
$ws.Range("M60").Value = -5839
$ws.Range("N60").Value = -11667

$ws.Range("H61").Value = 18500

$ws.Range("H113").Value = 2913.4092
$ws.Range("I113").Value = 2700.1
$ws.Range("J113").Value = 3091.1667
$ws.Range("K113").Value = 2700.1
$ws.Range("L113").Value = 3091.1667
$ws.Range("M113").Value = -530.0999999999999
$ws.Range("N113").Value = -7431.1667


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1499868.8
$ws.Range("I113").Value = 3448840.5
$ws.Range("J113").Value = 659.61536
$ws.Range("K113").Value = 10346521.5
$ws.Range("L113").Value = 1978.84608
$ws.Range("M113").Value = -10344351.5
$ws.Range("N113").Value = -6318.84608

$ws.Range("H132").Value = 1582.8
$ws.Range("J132").Value = 2003.3334
$ws.Range("L132").Value = 18030.0006
$ws.Range("N132").Value = -23090.0006


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 400.3846
$ws.Range("I107").Value = 363.81818
$ws.Range("J107").Value = 601.5
$ws.Range("K107").Value = 363.81818
$ws.Range("L107").Value = 601.5
$ws.Range("M107").Value = 1556.18182
$ws.Range("N107").Value = -4441.5

$ws.Range("H126").Value = 6671039
$ws.Range("I126").Value = 10004106
$ws.Range("J126").Value = 4905.6
$ws.Range("K126").Value = 30012318
$ws.Range("L126").Value = 14716.8
$ws.Range("M126").Value = -30009848
$ws.Range("N126").Value = -19656.8

$ws.Range("H132").Value = 3427.7273
$ws.Range("I132").Value = 4204.8
$ws.Range("J132").Value = 2780.1667
$ws.Range("K132").Value = 12614.4
$ws.Range("L132").Value = 8340.500100000001
$ws.Range("M132").Value = -10084.4
$ws.Range("N132").Value = -13400.5001


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1800
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 1800
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 3600
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -5722

$ws.Range("H84").Value = 1800
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 1800
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 18000
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -28608

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

